# Regenerate orders with updated distance/size codes.
# Rule (applied to every text cell in the used range):
#   D51 -> D55
#   D80 -> D86
#   D64 -> D69
#   S30 -> S31

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count

for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2

        if ($val -ne $null -and $val.GetType().Name -eq "String") {
            $newVal = $val
            $newVal = $newVal -replace "D51", "D55"
            $newVal = $newVal -replace "D80", "D86"
            $newVal = $newVal -replace "D64", "D69"
            $newVal = $newVal -replace "S30", "S31"

            if ($newVal -ne $val) {
                $cell.Value = $newVal
            }
        }
    }
}
